$p = $ppt.ActivePresentation

# ---------------------------------------------------------------------------
# 1) Three tables (slides 14, 15, 16) get their table style switched from
#    Table_0 ({BE94DCFC-8D0C-4CC5-89D2-1FC04648E48F}) to the built-in
#    "Medium Style 2 - Accent 1" gallery style ({E5698341-D5CF-4B34-97B0-D78AAAD0A895}).
# ---------------------------------------------------------------------------
$tableSlideIndexes = @(14, 15, 16)
$newStyleId = "{E5698341-D5CF-4B34-97B0-D78AAAD0A895}"

foreach ($slideIdx in $tableSlideIndexes) {
    $slide = $p.Slides.Item($slideIdx)
    for ($i = 1; $i -le $slide.Shapes.Count; $i++) {
        $shape = $slide.Shapes.Item($i)
        if ($shape.HasTable) {
            $shape.Table.ApplyStyle($newStyleId)
        }
    }
}

# ---------------------------------------------------------------------------
# 2) The deck's colour theme is switched from the custom "Integral" / "Red
#    Violet" scheme over to the standard "Office Theme" colour scheme.
#    (Font scheme / format scheme stay the untouched "Office" defaults they
#    already were.) Updating it through any slide's ThemeColorScheme updates
#    the one colour scheme backing the whole deck (master + every slide).
# ---------------------------------------------------------------------------
$officeColors = @{
    1  = 0         # dk1      000000
    2  = 16777215  # lt1      FFFFFF
    3  = 6968388   # dk2      44546A
    4  = 15132391  # lt2      E7E6E6
    5  = 13998939  # accent1  5B9BD5
    6  = 3243501   # accent2  ED7D31
    7  = 10855845  # accent3  A5A5A5
    8  = 49407     # accent4  FFC000
    9  = 12874308  # accent5  4472C4
    10 = 4697456   # accent6  70AD47
    11 = 12673797  # hlink    0563C1
    12 = 7491477   # folHlink 954F72
}

$tcs = $p.Slides.Item(1).ThemeColorScheme
foreach ($idx in $officeColors.Keys) {
    $tcs.Item($idx).RGB = $officeColors[$idx]
}
